$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  0.127881588408715,  0.04240448674262143, 0.8054896365839992, 8.660232485948974,   9.63600819768431),
    @(3,  3.230985683306322,  1.667794583268128,   0.8054896365839992, 8.660232485948974,   14.36450238910742),
    @(4,  3.230985683306322,  1.667794583268128,   0.8054896365839992, 0.496779210170732,   6.201049113329182),
    @(5,  3.230985683306322,  1.667794583268128,   0.1575252929769615, 0.496779210170732,   5.553084769722144),
    @(6,  0.127881588408715,  0.04240448674262143, 0.8054896365839992, 0.496779210170732,   1.472554921906068),
    @(7,  3.230985683306322,  10.29869402782916,   3.900430680208489,  8.660232485948974,   26.09034287729295),
    @(8,  3.230985683306322,  1.667794583268128,   0.1575252929769615, 0.496779210170732,   5.553084769722144),
    @(9,  3.230985683306322,  1.667794583268128,   0.8054896365839992, 0.496779210170732,   6.201049113329182),
    @(10, 1.459612070389937,  10.29869402782916,   3.900430680208489,  8.660232485948974,   24.31896926437656),
    @(11, 3.230985683306322,  1.667794583268128,   0.8054896365839992, 0.496779210170732,   6.201049113329182),
    @(12, 1.459612070389937,  1.667794583268128,   3.900430680208489,  0.496779210170732,   7.524616544037286),
    @(13, 0.04763786555579896,10.29869402782916,   3.900430680208489,  645.3272768299601,   659.5740394035535)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
